# LOM3255.xlsx update
# - Insert a new row above the old "Programa resumido:" row to hold a new
#   "Docentes responsáveis:" label (everything from that point down shifts
#   by one row).
# - Refresh several B/C value cells with new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 12 (shifts rows 12-20 down to 13-21) and label it.
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"

# Objetivos: now shows the responsible professor's id/name.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Programa resumido: now shows the term/periodicity.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: now shows the activation date.
$ws.Range("B15").Value = "15/07/2015"
$ws.Range("C15").Value = "15/07/2015"

# Método: now shows the responsible professor's id/name.
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
